$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = "'118"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'344175.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("C23").Value = "'260"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'960488.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("C26").Value = "'34"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'103572.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("C31").Value = "'48"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'110000.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("C32").Value = "'15"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'40500.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("C33").Value = "'27"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'78000.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("C35").Value = "'157"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'427908.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("C37").Value = "'330"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'1247535.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("C39").Value = "'5"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'24000.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("C40").Value = "'15"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'33653.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("C41").Value = "'40"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'123000.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("C43").Value = "'26"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'64971.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("C44").Value = "'4"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'8000.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("C46").Value = "'73"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'165893.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("C47").Value = "'13"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'46500.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("C59").Value = "'6"
$ws.Range("C59").Style = "Normal"
$ws.Range("D59").Value = "'12000.00"
$ws.Range("D59").Style = "Normal"
$ws.Range("C75").Value = "'30"
$ws.Range("C75").Style = "Normal"
$ws.Range("D75").Value = "'88500.00"
$ws.Range("D75").Style = "Normal"
$ws.Range("C76").Value = "'77"
$ws.Range("C76").Style = "Normal"
$ws.Range("D76").Value = "'204987.00"
$ws.Range("D76").Style = "Normal"
$ws.Range("C77").Value = "'175"
$ws.Range("C77").Style = "Normal"
$ws.Range("D77").Value = "'464722.00"
$ws.Range("D77").Style = "Normal"
$ws.Range("C79").Value = "'391"
$ws.Range("C79").Style = "Normal"
$ws.Range("D79").Value = "'1444130.70"
$ws.Range("D79").Style = "Normal"
$ws.Range("C83").Value = "'60"
$ws.Range("C83").Style = "Normal"
$ws.Range("D83").Value = "'195200.00"
$ws.Range("D83").Style = "Normal"
$ws.Range("C88").Value = "'92"
$ws.Range("C88").Style = "Normal"
$ws.Range("D88").Value = "'220500.00"
$ws.Range("D88").Style = "Normal"
$ws.Range("C121").Value = "'216"
$ws.Range("C121").Style = "Normal"
$ws.Range("D121").Value = "'590208.00"
$ws.Range("D121").Style = "Normal"
$ws.Range("C122").Value = "'49"
$ws.Range("C122").Style = "Normal"
$ws.Range("D122").Value = "'138827.58"
$ws.Range("D122").Style = "Normal"
$ws.Range("C123").Value = "'399"
$ws.Range("C123").Style = "Normal"
$ws.Range("D123").Value = "'1589448.95"
$ws.Range("D123").Style = "Normal"
$ws.Range("C131").Value = "'65"
$ws.Range("C131").Style = "Normal"
$ws.Range("D131").Value = "'260773.75"
$ws.Range("D131").Style = "Normal"
$ws.Range("C134").Value = "'190"
$ws.Range("C134").Style = "Normal"
$ws.Range("D134").Value = "'498120.00"
$ws.Range("D134").Style = "Normal"
$ws.Range("C136").Value = "'16"
$ws.Range("C136").Style = "Normal"
$ws.Range("D136").Value = "'36500.00"
$ws.Range("D136").Style = "Normal"
$ws.Range("C137").Value = "'507"
$ws.Range("C137").Style = "Normal"
$ws.Range("D137").Value = "'1258146.00"
$ws.Range("D137").Style = "Normal"
$ws.Range("C138").Value = "'1595"
$ws.Range("C138").Style = "Normal"
$ws.Range("D138").Value = "'4224883.03"
$ws.Range("D138").Style = "Normal"
$ws.Range("C139").Value = "'1847"
$ws.Range("C139").Style = "Normal"
$ws.Range("D139").Value = "'4468785.61"
$ws.Range("D139").Style = "Normal"
$ws.Range("C140").Value = "'2066"
$ws.Range("C140").Style = "Normal"
$ws.Range("D140").Value = "'7409219.66"
$ws.Range("D140").Style = "Normal"
$ws.Range("C141").Value = "'306"
$ws.Range("C141").Style = "Normal"
$ws.Range("D141").Value = "'805652.04"
$ws.Range("D141").Style = "Normal"
$ws.Range("C142").Value = "'112"
$ws.Range("C142").Style = "Normal"
$ws.Range("D142").Value = "'271000.00"
$ws.Range("D142").Style = "Normal"
$ws.Range("C143").Value = "'221"
$ws.Range("C143").Style = "Normal"
$ws.Range("D143").Value = "'543800.00"
$ws.Range("D143").Style = "Normal"
$ws.Range("C144").Value = "'890"
$ws.Range("C144").Style = "Normal"
$ws.Range("D144").Value = "'2258188.25"
$ws.Range("D144").Style = "Normal"
$ws.Range("C145").Value = "'424"
$ws.Range("C145").Style = "Normal"
$ws.Range("D145").Value = "'1142753.49"
$ws.Range("D145").Style = "Normal"
$ws.Range("C146").Value = "'321"
$ws.Range("C146").Style = "Normal"
$ws.Range("D146").Value = "'791700.16"
$ws.Range("D146").Style = "Normal"
$ws.Range("C147").Value = "'133"
$ws.Range("C147").Style = "Normal"
$ws.Range("D147").Value = "'317000.00"
$ws.Range("D147").Style = "Normal"
$ws.Range("C148").Value = "'315"
$ws.Range("C148").Style = "Normal"
$ws.Range("D148").Value = "'907421.68"
$ws.Range("D148").Style = "Normal"
$ws.Range("C149").Value = "'771"
$ws.Range("C149").Style = "Normal"
$ws.Range("D149").Value = "'1820612.82"
$ws.Range("D149").Style = "Normal"
$ws.Range("C167").Value = "'16"
$ws.Range("C167").Style = "Normal"
$ws.Range("D167").Value = "'47500.00"
$ws.Range("D167").Style = "Normal"
$ws.Range("C169").Value = "'155"
$ws.Range("C169").Style = "Normal"
$ws.Range("D169").Value = "'400000.00"
$ws.Range("D169").Style = "Normal"
$ws.Range("C171").Value = "'286"
$ws.Range("C171").Style = "Normal"
$ws.Range("D171").Value = "'986840.94"
$ws.Range("D171").Style = "Normal"
$ws.Range("C175").Value = "'50"
$ws.Range("C175").Style = "Normal"
$ws.Range("D175").Value = "'171058.79"
$ws.Range("D175").Style = "Normal"
$ws.Range("C176").Value = "'18"
$ws.Range("C176").Style = "Normal"
$ws.Range("D176").Value = "'50000.00"
$ws.Range("D176").Style = "Normal"
$ws.Range("C178").Value = "'11"
$ws.Range("C178").Style = "Normal"
$ws.Range("D178").Value = "'26500.00"
$ws.Range("D178").Style = "Normal"
$ws.Range("C179").Value = "'49"
$ws.Range("C179").Style = "Normal"
$ws.Range("D179").Value = "'182307.95"
$ws.Range("D179").Style = "Normal"
$ws.Range("C180").Value = "'65"
$ws.Range("C180").Style = "Normal"
$ws.Range("D180").Value = "'155000.00"
$ws.Range("D180").Style = "Normal"
$ws.Range("C183").Value = "'108"
$ws.Range("C183").Style = "Normal"
$ws.Range("D183").Value = "'285000.00"
$ws.Range("D183").Style = "Normal"
$ws.Range("C184").Value = "'333"
$ws.Range("C184").Style = "Normal"
$ws.Range("D184").Value = "'899788.00"
$ws.Range("D184").Style = "Normal"
$ws.Range("C185").Value = "'30"
$ws.Range("C185").Style = "Normal"
$ws.Range("D185").Value = "'91574.12"
$ws.Range("D185").Style = "Normal"
$ws.Range("C186").Value = "'582"
$ws.Range("C186").Style = "Normal"
$ws.Range("D186").Value = "'2027674.27"
$ws.Range("D186").Style = "Normal"
$ws.Range("C190").Value = "'148"
$ws.Range("C190").Style = "Normal"
$ws.Range("D190").Value = "'430633.00"
$ws.Range("D190").Style = "Normal"
$ws.Range("C194").Value = "'102"
$ws.Range("C194").Style = "Normal"
$ws.Range("D194").Value = "'426180.50"
$ws.Range("D194").Style = "Normal"
$ws.Range("C211").Value = "'17"
$ws.Range("C211").Style = "Normal"
$ws.Range("D211").Value = "'49078.00"
$ws.Range("D211").Style = "Normal"
$ws.Range("C214").Value = "'156"
$ws.Range("C214").Style = "Normal"
$ws.Range("D214").Value = "'427905.00"
$ws.Range("D214").Style = "Normal"
$ws.Range("C216").Value = "'308"
$ws.Range("C216").Style = "Normal"
$ws.Range("D216").Value = "'1017845.50"
$ws.Range("D216").Style = "Normal"
$ws.Range("C218").Value = "'6"
$ws.Range("C218").Style = "Normal"
$ws.Range("D218").Value = "'20000.00"
$ws.Range("D218").Style = "Normal"
$ws.Range("C220").Value = "'74"
$ws.Range("C220").Style = "Normal"
$ws.Range("D220").Value = "'207687.09"
$ws.Range("D220").Style = "Normal"
$ws.Range("C221").Value = "'25"
$ws.Range("C221").Style = "Normal"
$ws.Range("D221").Value = "'78587.00"
$ws.Range("D221").Style = "Normal"
$ws.Range("C223").Value = "'10"
$ws.Range("C223").Style = "Normal"
$ws.Range("D223").Value = "'29000.00"
$ws.Range("D223").Style = "Normal"
$ws.Range("C224").Value = "'48"
$ws.Range("C224").Style = "Normal"
$ws.Range("D224").Value = "'153270.00"
$ws.Range("D224").Style = "Normal"
$ws.Range("C229").Value = "'443"
$ws.Range("C229").Style = "Normal"
$ws.Range("D229").Value = "'1144583.00"
$ws.Range("D229").Style = "Normal"
